$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E235").NumberFormat = "@"

$ws.Cells.Item(1, 4).Value = 'codeforiati:category-code'
$ws.Cells.Item(1, 5).Value = 'codeforiati:group-code'
$ws.Cells.Item(1, 6).Value = 'codeforiati:group-name'
$ws.Cells.Item(1, 7).Value = 'codeforiati:category-name'
$ws.Cells.Item(2, 4).Value = '111'
$ws.Cells.Item(2, 5).Value = '110'
$ws.Cells.Item(2, 6).Value = 'Education'
$ws.Cells.Item(2, 7).Value = 'Education, Level Unspecified'
$ws.Cells.Item(3, 4).Value = '111'
$ws.Cells.Item(3, 5).Value = '110'
$ws.Cells.Item(3, 6).Value = 'Education'
$ws.Cells.Item(3, 7).Value = 'Education, Level Unspecified'
$ws.Cells.Item(4, 4).Value = '111'
$ws.Cells.Item(4, 5).Value = '110'
$ws.Cells.Item(4, 6).Value = 'Education'
$ws.Cells.Item(4, 7).Value = 'Education, Level Unspecified'
$ws.Cells.Item(5, 4).Value = '111'
$ws.Cells.Item(5, 5).Value = '110'
$ws.Cells.Item(5, 6).Value = 'Education'
$ws.Cells.Item(5, 7).Value = 'Education, Level Unspecified'
$ws.Cells.Item(6, 4).Value = '112'
$ws.Cells.Item(6, 5).Value = '110'
$ws.Cells.Item(6, 6).Value = 'Education'
$ws.Cells.Item(6, 7).Value = 'Basic Education'
$ws.Cells.Item(7, 4).Value = '112'
$ws.Cells.Item(7, 5).Value = '110'
$ws.Cells.Item(7, 6).Value = 'Education'
$ws.Cells.Item(7, 7).Value = 'Basic Education'
$ws.Cells.Item(8, 4).Value = '112'
$ws.Cells.Item(8, 5).Value = '110'
$ws.Cells.Item(8, 6).Value = 'Education'
$ws.Cells.Item(8, 7).Value = 'Basic Education'
$ws.Cells.Item(9, 4).Value = '112'
$ws.Cells.Item(9, 5).Value = '110'
$ws.Cells.Item(9, 6).Value = 'Education'
$ws.Cells.Item(9, 7).Value = 'Basic Education'
$ws.Cells.Item(10, 4).Value = '112'
$ws.Cells.Item(10, 5).Value = '110'
$ws.Cells.Item(10, 6).Value = 'Education'
$ws.Cells.Item(10, 7).Value = 'Basic Education'
$ws.Cells.Item(11, 4).Value = '112'
$ws.Cells.Item(11, 5).Value = '110'
$ws.Cells.Item(11, 6).Value = 'Education'
$ws.Cells.Item(11, 7).Value = 'Basic Education'
$ws.Cells.Item(12, 4).Value = '112'
$ws.Cells.Item(12, 5).Value = '110'
$ws.Cells.Item(12, 6).Value = 'Education'
$ws.Cells.Item(12, 7).Value = 'Basic Education'
$ws.Cells.Item(13, 4).Value = '113'
$ws.Cells.Item(13, 5).Value = '110'
$ws.Cells.Item(13, 6).Value = 'Education'
$ws.Cells.Item(13, 7).Value = 'Secondary Education'
$ws.Cells.Item(14, 4).Value = '113'
$ws.Cells.Item(14, 5).Value = '110'
$ws.Cells.Item(14, 6).Value = 'Education'
$ws.Cells.Item(14, 7).Value = 'Secondary Education'
$ws.Cells.Item(15, 4).Value = '114'
$ws.Cells.Item(15, 5).Value = '110'
$ws.Cells.Item(15, 6).Value = 'Education'
$ws.Cells.Item(15, 7).Value = 'Post-Secondary Education'
$ws.Cells.Item(16, 4).Value = '114'
$ws.Cells.Item(16, 5).Value = '110'
$ws.Cells.Item(16, 6).Value = 'Education'
$ws.Cells.Item(16, 7).Value = 'Post-Secondary Education'
$ws.Cells.Item(17, 4).Value = '121'
$ws.Cells.Item(17, 5).Value = '120'
$ws.Cells.Item(17, 6).Value = 'Health'
$ws.Cells.Item(17, 7).Value = 'Health, General'
$ws.Cells.Item(18, 4).Value = '121'
$ws.Cells.Item(18, 5).Value = '120'
$ws.Cells.Item(18, 6).Value = 'Health'
$ws.Cells.Item(18, 7).Value = 'Health, General'
$ws.Cells.Item(19, 4).Value = '121'
$ws.Cells.Item(19, 5).Value = '120'
$ws.Cells.Item(19, 6).Value = 'Health'
$ws.Cells.Item(19, 7).Value = 'Health, General'
$ws.Cells.Item(20, 4).Value = '121'
$ws.Cells.Item(20, 5).Value = '120'
$ws.Cells.Item(20, 6).Value = 'Health'
$ws.Cells.Item(20, 7).Value = 'Health, General'
$ws.Cells.Item(21, 4).Value = '122'
$ws.Cells.Item(21, 5).Value = '120'
$ws.Cells.Item(21, 6).Value = 'Health'
$ws.Cells.Item(21, 7).Value = 'Basic Health'
$ws.Cells.Item(22, 4).Value = '122'
$ws.Cells.Item(22, 5).Value = '120'
$ws.Cells.Item(22, 6).Value = 'Health'
$ws.Cells.Item(22, 7).Value = 'Basic Health'
$ws.Cells.Item(23, 4).Value = '122'
$ws.Cells.Item(23, 5).Value = '120'
$ws.Cells.Item(23, 6).Value = 'Health'
$ws.Cells.Item(23, 7).Value = 'Basic Health'
$ws.Cells.Item(24, 4).Value = '122'
$ws.Cells.Item(24, 5).Value = '120'
$ws.Cells.Item(24, 6).Value = 'Health'
$ws.Cells.Item(24, 7).Value = 'Basic Health'
$ws.Cells.Item(25, 4).Value = '122'
$ws.Cells.Item(25, 5).Value = '120'
$ws.Cells.Item(25, 6).Value = 'Health'
$ws.Cells.Item(25, 7).Value = 'Basic Health'
$ws.Cells.Item(26, 4).Value = '122'
$ws.Cells.Item(26, 5).Value = '120'
$ws.Cells.Item(26, 6).Value = 'Health'
$ws.Cells.Item(26, 7).Value = 'Basic Health'
$ws.Cells.Item(27, 4).Value = '122'
$ws.Cells.Item(27, 5).Value = '120'
$ws.Cells.Item(27, 6).Value = 'Health'
$ws.Cells.Item(27, 7).Value = 'Basic Health'
$ws.Cells.Item(28, 4).Value = '122'
$ws.Cells.Item(28, 5).Value = '120'
$ws.Cells.Item(28, 6).Value = 'Health'
$ws.Cells.Item(28, 7).Value = 'Basic Health'
$ws.Cells.Item(29, 4).Value = '122'
$ws.Cells.Item(29, 5).Value = '120'
$ws.Cells.Item(29, 6).Value = 'Health'
$ws.Cells.Item(29, 7).Value = 'Basic Health'
$ws.Cells.Item(30, 4).Value = '123'
$ws.Cells.Item(30, 5).Value = '120'
$ws.Cells.Item(30, 6).Value = 'Health'
$ws.Cells.Item(30, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(31, 4).Value = '123'
$ws.Cells.Item(31, 5).Value = '120'
$ws.Cells.Item(31, 6).Value = 'Health'
$ws.Cells.Item(31, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(32, 4).Value = '123'
$ws.Cells.Item(32, 5).Value = '120'
$ws.Cells.Item(32, 6).Value = 'Health'
$ws.Cells.Item(32, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(33, 4).Value = '123'
$ws.Cells.Item(33, 5).Value = '120'
$ws.Cells.Item(33, 6).Value = 'Health'
$ws.Cells.Item(33, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(34, 4).Value = '123'
$ws.Cells.Item(34, 5).Value = '120'
$ws.Cells.Item(34, 6).Value = 'Health'
$ws.Cells.Item(34, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(35, 4).Value = '123'
$ws.Cells.Item(35, 5).Value = '120'
$ws.Cells.Item(35, 6).Value = 'Health'
$ws.Cells.Item(35, 7).Value = 'Non-communicable diseases (NCDs)'
$ws.Cells.Item(36, 4).Value = '130'
$ws.Cells.Item(36, 5).Value = '130'
$ws.Cells.Item(36, 6).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(36, 7).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(37, 4).Value = '130'
$ws.Cells.Item(37, 5).Value = '130'
$ws.Cells.Item(37, 6).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(37, 7).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(38, 4).Value = '130'
$ws.Cells.Item(38, 5).Value = '130'
$ws.Cells.Item(38, 6).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(38, 7).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(39, 4).Value = '130'
$ws.Cells.Item(39, 5).Value = '130'
$ws.Cells.Item(39, 6).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(39, 7).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(40, 4).Value = '130'
$ws.Cells.Item(40, 5).Value = '130'
$ws.Cells.Item(40, 6).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(40, 7).Value = 'Population Policies/Programmes & Reproductive Health'
$ws.Cells.Item(41, 4).Value = '140'
$ws.Cells.Item(41, 5).Value = '140'
$ws.Cells.Item(41, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(41, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(42, 4).Value = '140'
$ws.Cells.Item(42, 5).Value = '140'
$ws.Cells.Item(42, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(42, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(43, 4).Value = '140'
$ws.Cells.Item(43, 5).Value = '140'
$ws.Cells.Item(43, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(43, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(44, 4).Value = '140'
$ws.Cells.Item(44, 5).Value = '140'
$ws.Cells.Item(44, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(44, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(45, 4).Value = '140'
$ws.Cells.Item(45, 5).Value = '140'
$ws.Cells.Item(45, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(45, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(46, 4).Value = '140'
$ws.Cells.Item(46, 5).Value = '140'
$ws.Cells.Item(46, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(46, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(47, 4).Value = '140'
$ws.Cells.Item(47, 5).Value = '140'
$ws.Cells.Item(47, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(47, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(48, 4).Value = '140'
$ws.Cells.Item(48, 5).Value = '140'
$ws.Cells.Item(48, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(48, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(49, 4).Value = '140'
$ws.Cells.Item(49, 5).Value = '140'
$ws.Cells.Item(49, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(49, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(50, 4).Value = '140'
$ws.Cells.Item(50, 5).Value = '140'
$ws.Cells.Item(50, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(50, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(51, 4).Value = '140'
$ws.Cells.Item(51, 5).Value = '140'
$ws.Cells.Item(51, 6).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(51, 7).Value = 'Water Supply & Sanitation'
$ws.Cells.Item(52, 4).Value = '151'
$ws.Cells.Item(52, 5).Value = '150'
$ws.Cells.Item(52, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(52, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(53, 4).Value = '151'
$ws.Cells.Item(53, 5).Value = '150'
$ws.Cells.Item(53, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(53, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(54, 4).Value = '151'
$ws.Cells.Item(54, 5).Value = '150'
$ws.Cells.Item(54, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(54, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(55, 4).Value = '151'
$ws.Cells.Item(55, 5).Value = '150'
$ws.Cells.Item(55, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(55, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(56, 4).Value = '151'
$ws.Cells.Item(56, 5).Value = '150'
$ws.Cells.Item(56, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(56, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(57, 4).Value = '151'
$ws.Cells.Item(57, 5).Value = '150'
$ws.Cells.Item(57, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(57, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(58, 4).Value = '151'
$ws.Cells.Item(58, 5).Value = '150'
$ws.Cells.Item(58, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(58, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(59, 4).Value = '151'
$ws.Cells.Item(59, 5).Value = '150'
$ws.Cells.Item(59, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(59, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(60, 4).Value = '151'
$ws.Cells.Item(60, 5).Value = '150'
$ws.Cells.Item(60, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(60, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(61, 4).Value = '151'
$ws.Cells.Item(61, 5).Value = '150'
$ws.Cells.Item(61, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(61, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(62, 4).Value = '151'
$ws.Cells.Item(62, 5).Value = '150'
$ws.Cells.Item(62, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(62, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(63, 4).Value = '151'
$ws.Cells.Item(63, 5).Value = '150'
$ws.Cells.Item(63, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(63, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(64, 4).Value = '151'
$ws.Cells.Item(64, 5).Value = '150'
$ws.Cells.Item(64, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(64, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(65, 4).Value = '151'
$ws.Cells.Item(65, 5).Value = '150'
$ws.Cells.Item(65, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(65, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(66, 4).Value = '151'
$ws.Cells.Item(66, 5).Value = '150'
$ws.Cells.Item(66, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(66, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(67, 4).Value = '151'
$ws.Cells.Item(67, 5).Value = '150'
$ws.Cells.Item(67, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(67, 7).Value = 'Government & Civil Society-general'
$ws.Cells.Item(68, 4).Value = '152'
$ws.Cells.Item(68, 5).Value = '150'
$ws.Cells.Item(68, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(68, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(69, 4).Value = '152'
$ws.Cells.Item(69, 5).Value = '150'
$ws.Cells.Item(69, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(69, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(70, 4).Value = '152'
$ws.Cells.Item(70, 5).Value = '150'
$ws.Cells.Item(70, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(70, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(71, 4).Value = '152'
$ws.Cells.Item(71, 5).Value = '150'
$ws.Cells.Item(71, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(71, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(72, 4).Value = '152'
$ws.Cells.Item(72, 5).Value = '150'
$ws.Cells.Item(72, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(72, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(73, 4).Value = '152'
$ws.Cells.Item(73, 5).Value = '150'
$ws.Cells.Item(73, 6).Value = 'Government & Civil Society'
$ws.Cells.Item(73, 7).Value = 'Conflict, Peace & Security'
$ws.Cells.Item(74, 4).Value = '160'
$ws.Cells.Item(74, 5).Value = '160'
$ws.Cells.Item(74, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(74, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(75, 4).Value = '160'
$ws.Cells.Item(75, 5).Value = '160'
$ws.Cells.Item(75, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(75, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(76, 4).Value = '160'
$ws.Cells.Item(76, 5).Value = '160'
$ws.Cells.Item(76, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(76, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(77, 4).Value = '160'
$ws.Cells.Item(77, 5).Value = '160'
$ws.Cells.Item(77, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(77, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(78, 4).Value = '160'
$ws.Cells.Item(78, 5).Value = '160'
$ws.Cells.Item(78, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(78, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(79, 4).Value = '160'
$ws.Cells.Item(79, 5).Value = '160'
$ws.Cells.Item(79, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(79, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(80, 4).Value = '160'
$ws.Cells.Item(80, 5).Value = '160'
$ws.Cells.Item(80, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(80, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(81, 4).Value = '160'
$ws.Cells.Item(81, 5).Value = '160'
$ws.Cells.Item(81, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(81, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(82, 4).Value = '160'
$ws.Cells.Item(82, 5).Value = '160'
$ws.Cells.Item(82, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(82, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(83, 4).Value = '160'
$ws.Cells.Item(83, 5).Value = '160'
$ws.Cells.Item(83, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(83, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(84, 4).Value = '160'
$ws.Cells.Item(84, 5).Value = '160'
$ws.Cells.Item(84, 6).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(84, 7).Value = 'Other Social Infrastructure & Services'
$ws.Cells.Item(85, 4).Value = '210'
$ws.Cells.Item(85, 5).Value = '210'
$ws.Cells.Item(85, 6).Value = 'Transport & Storage'
$ws.Cells.Item(85, 7).Value = 'Transport & Storage'
$ws.Cells.Item(86, 4).Value = '210'
$ws.Cells.Item(86, 5).Value = '210'
$ws.Cells.Item(86, 6).Value = 'Transport & Storage'
$ws.Cells.Item(86, 7).Value = 'Transport & Storage'
$ws.Cells.Item(87, 4).Value = '210'
$ws.Cells.Item(87, 5).Value = '210'
$ws.Cells.Item(87, 6).Value = 'Transport & Storage'
$ws.Cells.Item(87, 7).Value = 'Transport & Storage'
$ws.Cells.Item(88, 4).Value = '210'
$ws.Cells.Item(88, 5).Value = '210'
$ws.Cells.Item(88, 6).Value = 'Transport & Storage'
$ws.Cells.Item(88, 7).Value = 'Transport & Storage'
$ws.Cells.Item(89, 4).Value = '210'
$ws.Cells.Item(89, 5).Value = '210'
$ws.Cells.Item(89, 6).Value = 'Transport & Storage'
$ws.Cells.Item(89, 7).Value = 'Transport & Storage'
$ws.Cells.Item(90, 4).Value = '210'
$ws.Cells.Item(90, 5).Value = '210'
$ws.Cells.Item(90, 6).Value = 'Transport & Storage'
$ws.Cells.Item(90, 7).Value = 'Transport & Storage'
$ws.Cells.Item(91, 4).Value = '210'
$ws.Cells.Item(91, 5).Value = '210'
$ws.Cells.Item(91, 6).Value = 'Transport & Storage'
$ws.Cells.Item(91, 7).Value = 'Transport & Storage'
$ws.Cells.Item(92, 4).Value = '220'
$ws.Cells.Item(92, 5).Value = '220'
$ws.Cells.Item(92, 6).Value = 'Communications'
$ws.Cells.Item(92, 7).Value = 'Communications'
$ws.Cells.Item(93, 4).Value = '220'
$ws.Cells.Item(93, 5).Value = '220'
$ws.Cells.Item(93, 6).Value = 'Communications'
$ws.Cells.Item(93, 7).Value = 'Communications'
$ws.Cells.Item(94, 4).Value = '220'
$ws.Cells.Item(94, 5).Value = '220'
$ws.Cells.Item(94, 6).Value = 'Communications'
$ws.Cells.Item(94, 7).Value = 'Communications'
$ws.Cells.Item(95, 4).Value = '220'
$ws.Cells.Item(95, 5).Value = '220'
$ws.Cells.Item(95, 6).Value = 'Communications'
$ws.Cells.Item(95, 7).Value = 'Communications'
$ws.Cells.Item(96, 4).Value = '231'
$ws.Cells.Item(96, 5).Value = '230'
$ws.Cells.Item(96, 6).Value = 'Energy'
$ws.Cells.Item(96, 7).Value = 'Energy Policy'
$ws.Cells.Item(97, 4).Value = '231'
$ws.Cells.Item(97, 5).Value = '230'
$ws.Cells.Item(97, 6).Value = 'Energy'
$ws.Cells.Item(97, 7).Value = 'Energy Policy'
$ws.Cells.Item(98, 4).Value = '231'
$ws.Cells.Item(98, 5).Value = '230'
$ws.Cells.Item(98, 6).Value = 'Energy'
$ws.Cells.Item(98, 7).Value = 'Energy Policy'
$ws.Cells.Item(99, 4).Value = '231'
$ws.Cells.Item(99, 5).Value = '230'
$ws.Cells.Item(99, 6).Value = 'Energy'
$ws.Cells.Item(99, 7).Value = 'Energy Policy'
$ws.Cells.Item(100, 4).Value = '232'
$ws.Cells.Item(100, 5).Value = '230'
$ws.Cells.Item(100, 6).Value = 'Energy'
$ws.Cells.Item(100, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(101, 4).Value = '232'
$ws.Cells.Item(101, 5).Value = '230'
$ws.Cells.Item(101, 6).Value = 'Energy'
$ws.Cells.Item(101, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(102, 4).Value = '232'
$ws.Cells.Item(102, 5).Value = '230'
$ws.Cells.Item(102, 6).Value = 'Energy'
$ws.Cells.Item(102, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(103, 4).Value = '232'
$ws.Cells.Item(103, 5).Value = '230'
$ws.Cells.Item(103, 6).Value = 'Energy'
$ws.Cells.Item(103, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(104, 4).Value = '232'
$ws.Cells.Item(104, 5).Value = '230'
$ws.Cells.Item(104, 6).Value = 'Energy'
$ws.Cells.Item(104, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(105, 4).Value = '232'
$ws.Cells.Item(105, 5).Value = '230'
$ws.Cells.Item(105, 6).Value = 'Energy'
$ws.Cells.Item(105, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(106, 4).Value = '232'
$ws.Cells.Item(106, 5).Value = '230'
$ws.Cells.Item(106, 6).Value = 'Energy'
$ws.Cells.Item(106, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(107, 4).Value = '232'
$ws.Cells.Item(107, 5).Value = '230'
$ws.Cells.Item(107, 6).Value = 'Energy'
$ws.Cells.Item(107, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(108, 4).Value = '232'
$ws.Cells.Item(108, 5).Value = '230'
$ws.Cells.Item(108, 6).Value = 'Energy'
$ws.Cells.Item(108, 7).Value = 'Energy generation, renewable sources'
$ws.Cells.Item(109, 4).Value = '233'
$ws.Cells.Item(109, 5).Value = '230'
$ws.Cells.Item(109, 6).Value = 'Energy'
$ws.Cells.Item(109, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(110, 4).Value = '233'
$ws.Cells.Item(110, 5).Value = '230'
$ws.Cells.Item(110, 6).Value = 'Energy'
$ws.Cells.Item(110, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(111, 4).Value = '233'
$ws.Cells.Item(111, 5).Value = '230'
$ws.Cells.Item(111, 6).Value = 'Energy'
$ws.Cells.Item(111, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(112, 4).Value = '233'
$ws.Cells.Item(112, 5).Value = '230'
$ws.Cells.Item(112, 6).Value = 'Energy'
$ws.Cells.Item(112, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(113, 4).Value = '233'
$ws.Cells.Item(113, 5).Value = '230'
$ws.Cells.Item(113, 6).Value = 'Energy'
$ws.Cells.Item(113, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(114, 4).Value = '233'
$ws.Cells.Item(114, 5).Value = '230'
$ws.Cells.Item(114, 6).Value = 'Energy'
$ws.Cells.Item(114, 7).Value = 'Energy generation, non-renewable sources'
$ws.Cells.Item(115, 4).Value = '234'
$ws.Cells.Item(115, 5).Value = '230'
$ws.Cells.Item(115, 6).Value = 'Energy'
$ws.Cells.Item(115, 7).Value = 'Hybrid energy plants'
$ws.Cells.Item(116, 4).Value = '235'
$ws.Cells.Item(116, 5).Value = '230'
$ws.Cells.Item(116, 6).Value = 'Energy'
$ws.Cells.Item(116, 7).Value = 'Nuclear energy plants'
$ws.Cells.Item(117, 4).Value = '236'
$ws.Cells.Item(117, 5).Value = '230'
$ws.Cells.Item(117, 6).Value = 'Energy'
$ws.Cells.Item(117, 7).Value = 'Energy distribution'
$ws.Cells.Item(118, 4).Value = '236'
$ws.Cells.Item(118, 5).Value = '230'
$ws.Cells.Item(118, 6).Value = 'Energy'
$ws.Cells.Item(118, 7).Value = 'Energy distribution'
$ws.Cells.Item(119, 4).Value = '236'
$ws.Cells.Item(119, 5).Value = '230'
$ws.Cells.Item(119, 6).Value = 'Energy'
$ws.Cells.Item(119, 7).Value = 'Energy distribution'
$ws.Cells.Item(120, 4).Value = '236'
$ws.Cells.Item(120, 5).Value = '230'
$ws.Cells.Item(120, 6).Value = 'Energy'
$ws.Cells.Item(120, 7).Value = 'Energy distribution'
$ws.Cells.Item(121, 4).Value = '236'
$ws.Cells.Item(121, 5).Value = '230'
$ws.Cells.Item(121, 6).Value = 'Energy'
$ws.Cells.Item(121, 7).Value = 'Energy distribution'
$ws.Cells.Item(122, 4).Value = '236'
$ws.Cells.Item(122, 5).Value = '230'
$ws.Cells.Item(122, 6).Value = 'Energy'
$ws.Cells.Item(122, 7).Value = 'Energy distribution'
$ws.Cells.Item(123, 4).Value = '236'
$ws.Cells.Item(123, 5).Value = '230'
$ws.Cells.Item(123, 6).Value = 'Energy'
$ws.Cells.Item(123, 7).Value = 'Energy distribution'
$ws.Cells.Item(124, 4).Value = '240'
$ws.Cells.Item(124, 5).Value = '240'
$ws.Cells.Item(124, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(124, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(125, 4).Value = '240'
$ws.Cells.Item(125, 5).Value = '240'
$ws.Cells.Item(125, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(125, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(126, 4).Value = '240'
$ws.Cells.Item(126, 5).Value = '240'
$ws.Cells.Item(126, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(126, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(127, 4).Value = '240'
$ws.Cells.Item(127, 5).Value = '240'
$ws.Cells.Item(127, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(127, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(128, 4).Value = '240'
$ws.Cells.Item(128, 5).Value = '240'
$ws.Cells.Item(128, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(128, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(129, 4).Value = '240'
$ws.Cells.Item(129, 5).Value = '240'
$ws.Cells.Item(129, 6).Value = 'Banking & Financial Services'
$ws.Cells.Item(129, 7).Value = 'Banking & Financial Services'
$ws.Cells.Item(130, 4).Value = '250'
$ws.Cells.Item(130, 5).Value = '250'
$ws.Cells.Item(130, 6).Value = 'Business & Other Services'
$ws.Cells.Item(130, 7).Value = 'Business & Other Services'
$ws.Cells.Item(131, 4).Value = '250'
$ws.Cells.Item(131, 5).Value = '250'
$ws.Cells.Item(131, 6).Value = 'Business & Other Services'
$ws.Cells.Item(131, 7).Value = 'Business & Other Services'
$ws.Cells.Item(132, 4).Value = '250'
$ws.Cells.Item(132, 5).Value = '250'
$ws.Cells.Item(132, 6).Value = 'Business & Other Services'
$ws.Cells.Item(132, 7).Value = 'Business & Other Services'
$ws.Cells.Item(133, 4).Value = '250'
$ws.Cells.Item(133, 5).Value = '250'
$ws.Cells.Item(133, 6).Value = 'Business & Other Services'
$ws.Cells.Item(133, 7).Value = 'Business & Other Services'
$ws.Cells.Item(134, 4).Value = '311'
$ws.Cells.Item(134, 5).Value = '310'
$ws.Cells.Item(134, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(134, 7).Value = 'Agriculture'
$ws.Cells.Item(135, 4).Value = '311'
$ws.Cells.Item(135, 5).Value = '310'
$ws.Cells.Item(135, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(135, 7).Value = 'Agriculture'
$ws.Cells.Item(136, 4).Value = '311'
$ws.Cells.Item(136, 5).Value = '310'
$ws.Cells.Item(136, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(136, 7).Value = 'Agriculture'
$ws.Cells.Item(137, 4).Value = '311'
$ws.Cells.Item(137, 5).Value = '310'
$ws.Cells.Item(137, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(137, 7).Value = 'Agriculture'
$ws.Cells.Item(138, 4).Value = '311'
$ws.Cells.Item(138, 5).Value = '310'
$ws.Cells.Item(138, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(138, 7).Value = 'Agriculture'
$ws.Cells.Item(139, 4).Value = '311'
$ws.Cells.Item(139, 5).Value = '310'
$ws.Cells.Item(139, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(139, 7).Value = 'Agriculture'
$ws.Cells.Item(140, 4).Value = '311'
$ws.Cells.Item(140, 5).Value = '310'
$ws.Cells.Item(140, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(140, 7).Value = 'Agriculture'
$ws.Cells.Item(141, 4).Value = '311'
$ws.Cells.Item(141, 5).Value = '310'
$ws.Cells.Item(141, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(141, 7).Value = 'Agriculture'
$ws.Cells.Item(142, 4).Value = '311'
$ws.Cells.Item(142, 5).Value = '310'
$ws.Cells.Item(142, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(142, 7).Value = 'Agriculture'
$ws.Cells.Item(143, 4).Value = '311'
$ws.Cells.Item(143, 5).Value = '310'
$ws.Cells.Item(143, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(143, 7).Value = 'Agriculture'
$ws.Cells.Item(144, 4).Value = '311'
$ws.Cells.Item(144, 5).Value = '310'
$ws.Cells.Item(144, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(144, 7).Value = 'Agriculture'
$ws.Cells.Item(145, 4).Value = '311'
$ws.Cells.Item(145, 5).Value = '310'
$ws.Cells.Item(145, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(145, 7).Value = 'Agriculture'
$ws.Cells.Item(146, 4).Value = '311'
$ws.Cells.Item(146, 5).Value = '310'
$ws.Cells.Item(146, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(146, 7).Value = 'Agriculture'
$ws.Cells.Item(147, 4).Value = '311'
$ws.Cells.Item(147, 5).Value = '310'
$ws.Cells.Item(147, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(147, 7).Value = 'Agriculture'
$ws.Cells.Item(148, 4).Value = '311'
$ws.Cells.Item(148, 5).Value = '310'
$ws.Cells.Item(148, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(148, 7).Value = 'Agriculture'
$ws.Cells.Item(149, 4).Value = '311'
$ws.Cells.Item(149, 5).Value = '310'
$ws.Cells.Item(149, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(149, 7).Value = 'Agriculture'
$ws.Cells.Item(150, 4).Value = '311'
$ws.Cells.Item(150, 5).Value = '310'
$ws.Cells.Item(150, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(150, 7).Value = 'Agriculture'
$ws.Cells.Item(151, 4).Value = '311'
$ws.Cells.Item(151, 5).Value = '310'
$ws.Cells.Item(151, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(151, 7).Value = 'Agriculture'
$ws.Cells.Item(152, 4).Value = '312'
$ws.Cells.Item(152, 5).Value = '310'
$ws.Cells.Item(152, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(152, 7).Value = 'Forestry'
$ws.Cells.Item(153, 4).Value = '312'
$ws.Cells.Item(153, 5).Value = '310'
$ws.Cells.Item(153, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(153, 7).Value = 'Forestry'
$ws.Cells.Item(154, 4).Value = '312'
$ws.Cells.Item(154, 5).Value = '310'
$ws.Cells.Item(154, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(154, 7).Value = 'Forestry'
$ws.Cells.Item(155, 4).Value = '312'
$ws.Cells.Item(155, 5).Value = '310'
$ws.Cells.Item(155, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(155, 7).Value = 'Forestry'
$ws.Cells.Item(156, 4).Value = '312'
$ws.Cells.Item(156, 5).Value = '310'
$ws.Cells.Item(156, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(156, 7).Value = 'Forestry'
$ws.Cells.Item(157, 4).Value = '312'
$ws.Cells.Item(157, 5).Value = '310'
$ws.Cells.Item(157, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(157, 7).Value = 'Forestry'
$ws.Cells.Item(158, 4).Value = '313'
$ws.Cells.Item(158, 5).Value = '310'
$ws.Cells.Item(158, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(158, 7).Value = 'Fishing'
$ws.Cells.Item(159, 4).Value = '313'
$ws.Cells.Item(159, 5).Value = '310'
$ws.Cells.Item(159, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(159, 7).Value = 'Fishing'
$ws.Cells.Item(160, 4).Value = '313'
$ws.Cells.Item(160, 5).Value = '310'
$ws.Cells.Item(160, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(160, 7).Value = 'Fishing'
$ws.Cells.Item(161, 4).Value = '313'
$ws.Cells.Item(161, 5).Value = '310'
$ws.Cells.Item(161, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(161, 7).Value = 'Fishing'
$ws.Cells.Item(162, 4).Value = '313'
$ws.Cells.Item(162, 5).Value = '310'
$ws.Cells.Item(162, 6).Value = 'Agriculture, Forestry, Fishing'
$ws.Cells.Item(162, 7).Value = 'Fishing'
$ws.Cells.Item(163, 4).Value = '321'
$ws.Cells.Item(163, 5).Value = '320'
$ws.Cells.Item(163, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(163, 7).Value = 'Industry'
$ws.Cells.Item(164, 4).Value = '321'
$ws.Cells.Item(164, 5).Value = '320'
$ws.Cells.Item(164, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(164, 7).Value = 'Industry'
$ws.Cells.Item(165, 4).Value = '321'
$ws.Cells.Item(165, 5).Value = '320'
$ws.Cells.Item(165, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(165, 7).Value = 'Industry'
$ws.Cells.Item(166, 4).Value = '321'
$ws.Cells.Item(166, 5).Value = '320'
$ws.Cells.Item(166, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(166, 7).Value = 'Industry'
$ws.Cells.Item(167, 4).Value = '321'
$ws.Cells.Item(167, 5).Value = '320'
$ws.Cells.Item(167, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(167, 7).Value = 'Industry'
$ws.Cells.Item(168, 4).Value = '321'
$ws.Cells.Item(168, 5).Value = '320'
$ws.Cells.Item(168, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(168, 7).Value = 'Industry'
$ws.Cells.Item(169, 4).Value = '321'
$ws.Cells.Item(169, 5).Value = '320'
$ws.Cells.Item(169, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(169, 7).Value = 'Industry'
$ws.Cells.Item(170, 4).Value = '321'
$ws.Cells.Item(170, 5).Value = '320'
$ws.Cells.Item(170, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(170, 7).Value = 'Industry'
$ws.Cells.Item(171, 4).Value = '321'
$ws.Cells.Item(171, 5).Value = '320'
$ws.Cells.Item(171, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(171, 7).Value = 'Industry'
$ws.Cells.Item(172, 4).Value = '321'
$ws.Cells.Item(172, 5).Value = '320'
$ws.Cells.Item(172, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(172, 7).Value = 'Industry'
$ws.Cells.Item(173, 4).Value = '321'
$ws.Cells.Item(173, 5).Value = '320'
$ws.Cells.Item(173, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(173, 7).Value = 'Industry'
$ws.Cells.Item(174, 4).Value = '321'
$ws.Cells.Item(174, 5).Value = '320'
$ws.Cells.Item(174, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(174, 7).Value = 'Industry'
$ws.Cells.Item(175, 4).Value = '321'
$ws.Cells.Item(175, 5).Value = '320'
$ws.Cells.Item(175, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(175, 7).Value = 'Industry'
$ws.Cells.Item(176, 4).Value = '321'
$ws.Cells.Item(176, 5).Value = '320'
$ws.Cells.Item(176, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(176, 7).Value = 'Industry'
$ws.Cells.Item(177, 4).Value = '321'
$ws.Cells.Item(177, 5).Value = '320'
$ws.Cells.Item(177, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(177, 7).Value = 'Industry'
$ws.Cells.Item(178, 4).Value = '321'
$ws.Cells.Item(178, 5).Value = '320'
$ws.Cells.Item(178, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(178, 7).Value = 'Industry'
$ws.Cells.Item(179, 4).Value = '321'
$ws.Cells.Item(179, 5).Value = '320'
$ws.Cells.Item(179, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(179, 7).Value = 'Industry'
$ws.Cells.Item(180, 4).Value = '321'
$ws.Cells.Item(180, 5).Value = '320'
$ws.Cells.Item(180, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(180, 7).Value = 'Industry'
$ws.Cells.Item(181, 4).Value = '321'
$ws.Cells.Item(181, 5).Value = '320'
$ws.Cells.Item(181, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(181, 7).Value = 'Industry'
$ws.Cells.Item(182, 4).Value = '322'
$ws.Cells.Item(182, 5).Value = '320'
$ws.Cells.Item(182, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(182, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(183, 4).Value = '322'
$ws.Cells.Item(183, 5).Value = '320'
$ws.Cells.Item(183, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(183, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(184, 4).Value = '322'
$ws.Cells.Item(184, 5).Value = '320'
$ws.Cells.Item(184, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(184, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(185, 4).Value = '322'
$ws.Cells.Item(185, 5).Value = '320'
$ws.Cells.Item(185, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(185, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(186, 4).Value = '322'
$ws.Cells.Item(186, 5).Value = '320'
$ws.Cells.Item(186, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(186, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(187, 4).Value = '322'
$ws.Cells.Item(187, 5).Value = '320'
$ws.Cells.Item(187, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(187, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(188, 4).Value = '322'
$ws.Cells.Item(188, 5).Value = '320'
$ws.Cells.Item(188, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(188, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(189, 4).Value = '322'
$ws.Cells.Item(189, 5).Value = '320'
$ws.Cells.Item(189, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(189, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(190, 4).Value = '322'
$ws.Cells.Item(190, 5).Value = '320'
$ws.Cells.Item(190, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(190, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(191, 4).Value = '322'
$ws.Cells.Item(191, 5).Value = '320'
$ws.Cells.Item(191, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(191, 7).Value = 'Mineral Resources & Mining'
$ws.Cells.Item(192, 4).Value = '323'
$ws.Cells.Item(192, 5).Value = '320'
$ws.Cells.Item(192, 6).Value = 'Industry, Mining, Construction'
$ws.Cells.Item(192, 7).Value = 'Construction'
$ws.Cells.Item(193, 4).Value = '331'
$ws.Cells.Item(193, 5).Value = '331'
$ws.Cells.Item(193, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(193, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(194, 4).Value = '331'
$ws.Cells.Item(194, 5).Value = '331'
$ws.Cells.Item(194, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(194, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(195, 4).Value = '331'
$ws.Cells.Item(195, 5).Value = '331'
$ws.Cells.Item(195, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(195, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(196, 4).Value = '331'
$ws.Cells.Item(196, 5).Value = '331'
$ws.Cells.Item(196, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(196, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(197, 4).Value = '331'
$ws.Cells.Item(197, 5).Value = '331'
$ws.Cells.Item(197, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(197, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(198, 4).Value = '331'
$ws.Cells.Item(198, 5).Value = '331'
$ws.Cells.Item(198, 6).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(198, 7).Value = 'Trade Policies & Regulations'
$ws.Cells.Item(199, 4).Value = '332'
$ws.Cells.Item(199, 5).Value = '332'
$ws.Cells.Item(199, 6).Value = 'Tourism'
$ws.Cells.Item(199, 7).Value = 'Tourism'
$ws.Cells.Item(200, 4).Value = '410'
$ws.Cells.Item(200, 5).Value = '410'
$ws.Cells.Item(200, 6).Value = 'General Environment Protection'
$ws.Cells.Item(200, 7).Value = 'General Environment Protection'
$ws.Cells.Item(201, 4).Value = '410'
$ws.Cells.Item(201, 5).Value = '410'
$ws.Cells.Item(201, 6).Value = 'General Environment Protection'
$ws.Cells.Item(201, 7).Value = 'General Environment Protection'
$ws.Cells.Item(202, 4).Value = '410'
$ws.Cells.Item(202, 5).Value = '410'
$ws.Cells.Item(202, 6).Value = 'General Environment Protection'
$ws.Cells.Item(202, 7).Value = 'General Environment Protection'
$ws.Cells.Item(203, 4).Value = '410'
$ws.Cells.Item(203, 5).Value = '410'
$ws.Cells.Item(203, 6).Value = 'General Environment Protection'
$ws.Cells.Item(203, 7).Value = 'General Environment Protection'
$ws.Cells.Item(204, 4).Value = '410'
$ws.Cells.Item(204, 5).Value = '410'
$ws.Cells.Item(204, 6).Value = 'General Environment Protection'
$ws.Cells.Item(204, 7).Value = 'General Environment Protection'
$ws.Cells.Item(205, 4).Value = '410'
$ws.Cells.Item(205, 5).Value = '410'
$ws.Cells.Item(205, 6).Value = 'General Environment Protection'
$ws.Cells.Item(205, 7).Value = 'General Environment Protection'
$ws.Cells.Item(206, 4).Value = '430'
$ws.Cells.Item(206, 5).Value = '430'
$ws.Cells.Item(206, 6).Value = 'Other Multisector'
$ws.Cells.Item(206, 7).Value = 'Other Multisector'
$ws.Cells.Item(207, 4).Value = '430'
$ws.Cells.Item(207, 5).Value = '430'
$ws.Cells.Item(207, 6).Value = 'Other Multisector'
$ws.Cells.Item(207, 7).Value = 'Other Multisector'
$ws.Cells.Item(208, 4).Value = '430'
$ws.Cells.Item(208, 5).Value = '430'
$ws.Cells.Item(208, 6).Value = 'Other Multisector'
$ws.Cells.Item(208, 7).Value = 'Other Multisector'
$ws.Cells.Item(209, 4).Value = '430'
$ws.Cells.Item(209, 5).Value = '430'
$ws.Cells.Item(209, 6).Value = 'Other Multisector'
$ws.Cells.Item(209, 7).Value = 'Other Multisector'
$ws.Cells.Item(210, 4).Value = '430'
$ws.Cells.Item(210, 5).Value = '430'
$ws.Cells.Item(210, 6).Value = 'Other Multisector'
$ws.Cells.Item(210, 7).Value = 'Other Multisector'
$ws.Cells.Item(211, 4).Value = '430'
$ws.Cells.Item(211, 5).Value = '430'
$ws.Cells.Item(211, 6).Value = 'Other Multisector'
$ws.Cells.Item(211, 7).Value = 'Other Multisector'
$ws.Cells.Item(212, 4).Value = '430'
$ws.Cells.Item(212, 5).Value = '430'
$ws.Cells.Item(212, 6).Value = 'Other Multisector'
$ws.Cells.Item(212, 7).Value = 'Other Multisector'
$ws.Cells.Item(213, 4).Value = '430'
$ws.Cells.Item(213, 5).Value = '430'
$ws.Cells.Item(213, 6).Value = 'Other Multisector'
$ws.Cells.Item(213, 7).Value = 'Other Multisector'
$ws.Cells.Item(214, 4).Value = '430'
$ws.Cells.Item(214, 5).Value = '430'
$ws.Cells.Item(214, 6).Value = 'Other Multisector'
$ws.Cells.Item(214, 7).Value = 'Other Multisector'
$ws.Cells.Item(215, 4).Value = '430'
$ws.Cells.Item(215, 5).Value = '430'
$ws.Cells.Item(215, 6).Value = 'Other Multisector'
$ws.Cells.Item(215, 7).Value = 'Other Multisector'
$ws.Cells.Item(216, 4).Value = '510'
$ws.Cells.Item(216, 5).Value = '510'
$ws.Cells.Item(216, 6).Value = 'General Budget Support'
$ws.Cells.Item(216, 7).Value = 'General Budget Support'
$ws.Cells.Item(217, 4).Value = '520'
$ws.Cells.Item(217, 5).Value = '520'
$ws.Cells.Item(217, 6).Value = 'Development Food Assistance'
$ws.Cells.Item(217, 7).Value = 'Development Food Assistance'
$ws.Cells.Item(218, 4).Value = '530'
$ws.Cells.Item(218, 5).Value = '530'
$ws.Cells.Item(218, 6).Value = 'Other Commodity Assistance'
$ws.Cells.Item(218, 7).Value = 'Other Commodity Assistance'
$ws.Cells.Item(219, 4).Value = '530'
$ws.Cells.Item(219, 5).Value = '530'
$ws.Cells.Item(219, 6).Value = 'Other Commodity Assistance'
$ws.Cells.Item(219, 7).Value = 'Other Commodity Assistance'
$ws.Cells.Item(220, 4).Value = '600'
$ws.Cells.Item(220, 5).Value = '600'
$ws.Cells.Item(220, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(220, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(221, 4).Value = '600'
$ws.Cells.Item(221, 5).Value = '600'
$ws.Cells.Item(221, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(221, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(222, 4).Value = '600'
$ws.Cells.Item(222, 5).Value = '600'
$ws.Cells.Item(222, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(222, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(223, 4).Value = '600'
$ws.Cells.Item(223, 5).Value = '600'
$ws.Cells.Item(223, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(223, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(224, 4).Value = '600'
$ws.Cells.Item(224, 5).Value = '600'
$ws.Cells.Item(224, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(224, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(225, 4).Value = '600'
$ws.Cells.Item(225, 5).Value = '600'
$ws.Cells.Item(225, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(225, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(226, 4).Value = '600'
$ws.Cells.Item(226, 5).Value = '600'
$ws.Cells.Item(226, 6).Value = 'Action Relating to Debt'
$ws.Cells.Item(226, 7).Value = 'Action Relating to Debt'
$ws.Cells.Item(227, 4).Value = '720'
$ws.Cells.Item(227, 5).Value = '720'
$ws.Cells.Item(227, 6).Value = 'Emergency Response'
$ws.Cells.Item(227, 7).Value = 'Emergency Response'
$ws.Cells.Item(228, 4).Value = '720'
$ws.Cells.Item(228, 5).Value = '720'
$ws.Cells.Item(228, 6).Value = 'Emergency Response'
$ws.Cells.Item(228, 7).Value = 'Emergency Response'
$ws.Cells.Item(229, 4).Value = '720'
$ws.Cells.Item(229, 5).Value = '720'
$ws.Cells.Item(229, 6).Value = 'Emergency Response'
$ws.Cells.Item(229, 7).Value = 'Emergency Response'
$ws.Cells.Item(230, 4).Value = '730'
$ws.Cells.Item(230, 5).Value = '730'
$ws.Cells.Item(230, 6).Value = 'Reconstruction Relief & Rehabilitation'
$ws.Cells.Item(230, 7).Value = 'Reconstruction Relief & Rehabilitation'
$ws.Cells.Item(231, 4).Value = '740'
$ws.Cells.Item(231, 5).Value = '740'
$ws.Cells.Item(231, 6).Value = 'Disaster Prevention & Preparedness'
$ws.Cells.Item(231, 7).Value = 'Disaster Prevention & Preparedness'
$ws.Cells.Item(232, 4).Value = '910'
$ws.Cells.Item(232, 5).Value = '910'
$ws.Cells.Item(232, 6).Value = 'Administrative Costs of Donors'
$ws.Cells.Item(232, 7).Value = 'Administrative Costs of Donors'
$ws.Cells.Item(233, 4).Value = '930'
$ws.Cells.Item(233, 5).Value = '930'
$ws.Cells.Item(233, 6).Value = 'Refugees in Donor Countries'
$ws.Cells.Item(233, 7).Value = 'Refugees in Donor Countries'
$ws.Cells.Item(234, 4).Value = '998'
$ws.Cells.Item(234, 5).Value = '998'
$ws.Cells.Item(234, 6).Value = 'Unallocated / Unspecified'
$ws.Cells.Item(234, 7).Value = 'Unallocated / Unspecified'
$ws.Cells.Item(235, 4).Value = '998'
$ws.Cells.Item(235, 5).Value = '998'
$ws.Cells.Item(235, 6).Value = 'Unallocated / Unspecified'
$ws.Cells.Item(235, 7).Value = 'Unallocated / Unspecified'
